# FedExShipments_PreProd.xlsx — "Changes of webdriver exception"
#
# The FedEx tracking numbers in column P (rows 2-26) are refreshed to a new
# batch of shipment tracking IDs. Each tracking number is stored as text
# (shared string), so when we write the replacement values — which happen
# to look like plain integers — we force the cell to stay text (NumberFormat
# "@") while we update its contents, then clear the explicit formatting
# again so the cell ends up exactly like the original: default style, text
# value, no number formatting applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTrackingNumbers = @(
    "320017962708",
    "320017962719",
    "320017962741",
    "320017962774",
    "320017962811",
    "320017962833",
    "320017962866",
    "320017962888",
    "320017962936",
    "320017962958",
    "320017962991",
    "320017963016",
    "320017963049",
    "320017963060",
    "320017963093",
    "320017963119",
    "320017963152",
    "320017963174",
    "320017963200",
    "320017963222",
    "320017963255",
    "320017963266",
    "320017963288",
    "320017963299",
    "320017963314"
)

$firstRow = 2
$col = 16  # column P

for ($i = 0; $i -lt $newTrackingNumbers.Length; $i++) {
    $row = $firstRow + $i
    $cell = $ws.Cells.Item($row, $col)

    # Force text storage so the numeric-looking tracking number isn't
    # reinterpreted as a number, then restore the default (no explicit)
    # cell style so formatting matches the surrounding cells.
    $cell.NumberFormat = "@"
    $cell.Characters().Text = $newTrackingNumbers[$i]
    $cell.ClearFormats()
}
